$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 3000
$ws.Range("I9").Value = 3000
$ws.Range("K9").Value = 3000
$ws.Range("M9").Value = -2831
$ws.Range("H33").Value = 697.9524
$ws.Range("I33").Value = 794.7059
$ws.Range("K33").Value = 794.7059
$ws.Range("M33").Value = -565.7059
$ws.Range("H53").Value = 175.76471
$ws.Range("J53").Value = 87.666664
$ws.Range("L53").Value = 87.666664
$ws.Range("N53").Value = -1361.666664
$ws.Range("H80").Value = 28570.445
$ws.Range("J80").Value = 45833.11
$ws.Range("L80").Value = 137499.33
$ws.Range("N80").Value = -139495.33
$ws.Range("H83").Value = 28570.445
$ws.Range("J83").Value = 45833.11
$ws.Range("L83").Value = 412497.99
$ws.Range("N83").Value = -422481.99
$ws.Range("H132").Value = 2226.1035
$ws.Range("I132").Value = 2252.0356
$ws.Range("K132").Value = 6756.1068
$ws.Range("M132").Value = -4226.1068
$ws.Range("H137").Value = 2975.6428
$ws.Range("I137").Value = 3239.889
$ws.Range("J137").Value = 2500
$ws.Range("K137").Value = 9719.667000000001
$ws.Range("L137").Value = 7500
$ws.Range("M137").Value = -7169.667000000001
$ws.Range("N137").Value = -12600

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 7238.0884
$ws.Range("I61").Value = 3760.28
$ws.Range("J61").Value = 16898.666
$ws.Range("K61").Value = 3760.28
$ws.Range("L61").Value = 16898.666
$ws.Range("M61").Value = -3548.28
$ws.Range("N61").Value = -17322.666
$ws.Range("H74").Value = 3091.762
$ws.Range("I74").Value = 1303.3077
$ws.Range("K74").Value = 1303.3077
$ws.Range("M74").Value = -429.3077000000001
$ws.Range("H77").Value = 3091.762
$ws.Range("I77").Value = 1303.3077
$ws.Range("K77").Value = 6516.538500000001
$ws.Range("M77").Value = -2148.538500000001
$ws.Range("H122").Value = 10609.107
$ws.Range("I122").Value = 12532.3
$ws.Range("K122").Value = 37596.89999999999
$ws.Range("M122").Value = -35146.89999999999
$ws.Range("H132").Value = 8872.843999999999
$ws.Range("I132").Value = 8740.866
$ws.Range("K132").Value = 26222.598
$ws.Range("M132").Value = -23692.598
$ws.Range("H136").Value = 7238.0884
$ws.Range("I136").Value = 3760.28
$ws.Range("J136").Value = 16898.666
$ws.Range("K136").Value = 11280.84
$ws.Range("L136").Value = 50695.99800000001
$ws.Range("M136").Value = -8730.84
$ws.Range("N136").Value = -55795.99800000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 29412008
$ws.Range("I80").Value = 83333660
$ws.Range("J80").Value = 194.18182
$ws.Range("K80").Value = 83333660
$ws.Range("L80").Value = 194.18182
$ws.Range("M80").Value = -83332662
$ws.Range("N80").Value = -2190.18182
$ws.Range("H83").Value = 29412008
$ws.Range("I83").Value = 83333660
$ws.Range("J83").Value = 194.18182
$ws.Range("K83").Value = 416668300
$ws.Range("L83").Value = 970.9091
$ws.Range("M83").Value = -416663308
$ws.Range("N83").Value = -10954.9091

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 4470.9
$ws.Range("I16").Value = 2510.7273
$ws.Range("J16").Value = 6866.6665
$ws.Range("K16").Value = 2510.7273
$ws.Range("L16").Value = 6866.6665
$ws.Range("M16").Value = -2223.7273
$ws.Range("N16").Value = -7440.6665
$ws.Range("H22").Value = 521.55554
$ws.Range("I22").Value = 350
$ws.Range("J22").Value = 736
$ws.Range("K22").Value = 350
$ws.Range("L22").Value = 736
$ws.Range("M22").Value = 0
$ws.Range("N22").Value = -1436
$ws.Range("H31").Value = 8032.9585
$ws.Range("I31").Value = 4093.923
$ws.Range("K31").Value = 4093.923
$ws.Range("M31").Value = -3798.923
$ws.Range("H34").Value = 8032.9585
$ws.Range("I34").Value = 4093.923
$ws.Range("K34").Value = 4093.923
$ws.Range("M34").Value = -3891.923
$ws.Range("H58").Value = 13520502
$ws.Range("I58").Value = 35716572
$ws.Range("J58").Value = 9849.261
$ws.Range("K58").Value = 35716572
$ws.Range("L58").Value = 9849.261
$ws.Range("M58").Value = -35716369
$ws.Range("N58").Value = -10255.261
$ws.Range("H107").Value = 1499.3889
$ws.Range("I107").Value = 1138.8518
$ws.Range("J107").Value = 2581
$ws.Range("K107").Value = 1138.8518
$ws.Range("L107").Value = 2581
$ws.Range("M107").Value = 781.1482000000001
$ws.Range("N107").Value = -6421
$ws.Range("H113").Value = 4470.9
$ws.Range("I113").Value = 2510.7273
$ws.Range("J113").Value = 6866.6665
$ws.Range("K113").Value = 2510.7273
$ws.Range("L113").Value = 6866.6665
$ws.Range("M113").Value = -340.7273
$ws.Range("N113").Value = -11206.6665
$ws.Range("H122").Value = 3255.95
$ws.Range("I122").Value = 3017.2
$ws.Range("J122").Value = 3972.2
$ws.Range("K122").Value = 9051.599999999999
$ws.Range("L122").Value = 11916.6
$ws.Range("M122").Value = -6601.599999999999
$ws.Range("N122").Value = -16816.6
$ws.Range("H136").Value = 13520502
$ws.Range("I136").Value = 35716572
$ws.Range("J136").Value = 9849.261
$ws.Range("K136").Value = 107149716
$ws.Range("L136").Value = 29547.783
$ws.Range("M136").Value = -107147166
$ws.Range("N136").Value = -34647.783

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 123939580
$ws.Range("I4").Value = 151481660
$ws.Range("K4").Value = 454444980
$ws.Range("M4").Value = -454444868
$ws.Range("H122").Value = 1665892.4
$ws.Range("I122").Value = 3537472.5
$ws.Range("J122").Value = 2265.5557
$ws.Range("K122").Value = 31837252.5
$ws.Range("L122").Value = 20390.0013
$ws.Range("M122").Value = -31834802.5
$ws.Range("N122").Value = -25290.0013

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H58").Value = 68428.2
$ws.Range("J58").Value = 68428.2
$ws.Range("L58").Value = 68428.2
$ws.Range("N58").Value = -68982.2
$ws.Range("H70").Value = 7997.706
$ws.Range("J70").Value = 8565.666999999999
$ws.Range("L70").Value = 8565.666999999999
$ws.Range("N70").Value = -9105.666999999999
$ws.Range("H73").Value = 7997.706
$ws.Range("J73").Value = 8565.666999999999
$ws.Range("L73").Value = 8565.666999999999
$ws.Range("N73").Value = -10437.667
$ws.Range("H113").Value = 7574.7354
$ws.Range("I113").Value = 4712.4165
$ws.Range("J113").Value = 9136
$ws.Range("K113").Value = 4712.4165
$ws.Range("L113").Value = 9136
$ws.Range("M113").Value = -2542.4165
$ws.Range("N113").Value = -13476
$ws.Range("H126").Value = 3442.5789
$ws.Range("I126").Value = 1413.4445
$ws.Range("J126").Value = 5268.8
$ws.Range("K126").Value = 4240.333500000001
$ws.Range("L126").Value = 15806.4
$ws.Range("M126").Value = -1770.333500000001
$ws.Range("N126").Value = -20746.4
$ws.Range("H132").Value = 4185
$ws.Range("I132").Value = 1775.9474
$ws.Range("J132").Value = 7999.3335
$ws.Range("K132").Value = 5327.8422
$ws.Range("L132").Value = 23998.0005
$ws.Range("M132").Value = -2797.8422
$ws.Range("N132").Value = -29058.0005

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5297.8887
$ws.Range("I7").Value = 3396.0908
$ws.Range("J7").Value = 8286.429
$ws.Range("K7").Value = 3396.0908
$ws.Range("L7").Value = 8286.429
$ws.Range("M7").Value = -3284.0908
$ws.Range("N7").Value = -8510.429
$ws.Range("H40").Value = 5387.185
$ws.Range("I40").Value = 4845.609
$ws.Range("K40").Value = 4845.609
$ws.Range("M40").Value = -4709.609
$ws.Range("H122").Value = 4263.8
$ws.Range("J122").Value = 6667.3335
$ws.Range("L122").Value = 20002.0005
$ws.Range("N122").Value = -24902.0005
$ws.Range("H126").Value = 5297.8887
$ws.Range("I126").Value = 3396.0908
$ws.Range("J126").Value = 8286.429
$ws.Range("K126").Value = 10188.2724
$ws.Range("L126").Value = 24859.287
$ws.Range("M126").Value = -7718.2724
$ws.Range("N126").Value = -29799.287
$ws.Range("H132").Value = 13895737
$ws.Range("I132").Value = 25003060
$ws.Range("J132").Value = 11583.3125
$ws.Range("K132").Value = 75009180
$ws.Range("L132").Value = 34749.9375
$ws.Range("M132").Value = -75006650
$ws.Range("N132").Value = -39809.9375

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 18188212
$ws.Range("J81").Value = 66686776
$ws.Range("L81").Value = 133373552
$ws.Range("N81").Value = -133375674
$ws.Range("H84").Value = 18188212
$ws.Range("J84").Value = 66686776
$ws.Range("L84").Value = 666867760
$ws.Range("N84").Value = -666878368
$ws.Range("H100").Value = 412
$ws.Range("I100").Value = 374.4
$ws.Range("K100").Value = 748.8
$ws.Range("M100").Value = -207.8
$ws.Range("H122").Value = 26530906
$ws.Range("I122").Value = 42004100
$ws.Range("J122").Value = 5428.5713
$ws.Range("K122").Value = 126012300
$ws.Range("L122").Value = 16285.7139
$ws.Range("M122").Value = -126009850
$ws.Range("N122").Value = -21185.7139
$ws.Range("H126").Value = 2363.2632
$ws.Range("I126").Value = 1659.4073
$ws.Range("J126").Value = 4090.9092
$ws.Range("K126").Value = 4978.2219
$ws.Range("L126").Value = 12272.7276
$ws.Range("M126").Value = -2508.2219
$ws.Range("N126").Value = -17212.7276
$ws.Range("H132").Value = 21767434
$ws.Range("I132").Value = 25007048
$ws.Range("J132").Value = 169999.67
$ws.Range("K132").Value = 75021144
$ws.Range("L132").Value = 509999.01
$ws.Range("M132").Value = -75018614
$ws.Range("N132").Value = -515059.01
